$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.588.37"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "2.071.43"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'231.75"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'57.89"
$ws.Range("E8").Value = "  -2.29%  "
$ws.Range("D9").Value = "'0.388"
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("D10").Value = "'0.0775"
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "2.375.88"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").Value = "'14.76"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "'21.19"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "'0.764"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").Value = "'5.32"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "2.075.17"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "37.515.86"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "'6.15"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'69.93"
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("D21").Value = "0.0₃0826"
$ws.Range("E21").Value = "  -2.92%  "
$ws.Range("D22").Value = "'227.08"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("E25").Value = "  -3.20%  "
$ws.Range("E26").Value = "  +3.62%  "
$ws.Range("D27").Value = "'169.41"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("E28").Value = "  -5.15%  "
$ws.Range("D29").Value = "'19.34"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("D30").Value = "'1.36"
$ws.Range("E30").Value = "  -4.79%  "
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("E32").Value = "  -3.58%  "
$ws.Range("D33").Value = "'0.0626"
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("D34").Value = "'4.64"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").Value = "'2.53"
$ws.Range("E35").Value = "  +1.09%  "
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("E37").Value = "  -4.03%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "'5.32"
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("D40").Value = "'0.0227"
$ws.Range("E40").Value = "  +3.48%  "
$ws.Range("D41").Value = "'98.23"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").Value = "'0.0957"
$ws.Range("E42").Value = "  -2.76%  "
$ws.Range("D43").Value = "1.489.56"
$ws.Range("E43").Value = "  +2.69%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D46").Value = "'16.61"
$ws.Range("E46").Value = "  -3.67%  "
$ws.Range("E47").Value = "  -2.62%  "
$ws.Range("D48").Value = "'4.01"
$ws.Range("E48").Value = "  -3.59%  "
$ws.Range("D49").Value = "'7.25"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("D51").Value = "2.261.91"
$ws.Range("E51").Value = "  -0.70%  "
